$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update horizontal/vertical timing inputs (Sync pulse & Back porch line counts)
$ws.Range("B16").Value = 12
$ws.Range("B18").Value = 33

# Recalculate dependent formulas
$excel.Calculate()

# Move the active selection to B16, as left by the editor
$ws.Activate()
$ws.Range("B16").Select()
